# إضافة حدث جديد في Card6 by admin at 2026-01-20 16:04:46
#
# The "Card6" log sheet keeps a small lookup table (rows 2-12, columns
# A-K) followed by a service/event log (columns L-O: Date / Event /
# Correction / Serviced by). Appending a new event adds one more row at
# the bottom of the sheet and normalises any still-empty data columns
# (B:K) on the row(s) involved to the literal text "nan", matching how
# the sheet has always represented "no value" in those columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card6")

# Locate the current last populated row (column A always holds the card
# number, "6", for every row in this sheet).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# The row that used to be last still had untouched/empty B:K cells;
# normalise them to "nan" the same way the rest of the sheet does.
$ws.Range("B" + $lastRow + ":K" + $lastRow).Value = "nan"

# New event row: card number in column A (kept as text, like the rest
# of the column), B:K left blank (no tone-range data on a log row),
# and the event details in L:O.
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("A" + $newRow).Value = "6"
$ws.Range("A" + $newRow).Style = "Normal"

$ws.Range("L" + $newRow).Value = "20/1/2026"
$ws.Range("M" + $newRow).Value = "زياره توكيل"
$ws.Range("N" + $newRow).Value = "تم تغير سوفت كرد لbc"
$ws.Range("O" + $newRow).Value = "م. احمد علي توكيل"
